$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G6").Value  = "Chikkamagaluru (Chikmagalur)"
$ws.Range("G12").Value = "Chikkamagaluru (Chikmagalur)"
$ws.Range("G13").Value = "Shivamogga (Shimoga)"
$ws.Range("G19").Value = "Chikkamagaluru (Chikmagalur)"
$ws.Range("G21").Value = "Bidar"
$ws.Range("G23").Value = "Shivamogga (Shimoga)"
$ws.Range("G26").Value = "Chikkamagaluru (Chikmagalur)"
$ws.Range("G28").Value = "Shivamogga (Shimoga)"
$ws.Range("G30").Value = "Vijayapura (Bijapur)"
$ws.Range("G31").Value = "Shivamogga (Shimoga)"
$ws.Range("G32").Value = "Chikkamagaluru (Chikmagalur)"
$ws.Range("G37").Value = "Shivamogga (Shimoga)"
$ws.Range("G38").Value = "Shivamogga (Shimoga)"
$ws.Range("G40").Value = "Chikkamagaluru (Chikmagalur)"
$ws.Range("G42").Value = "Chikkamagaluru (Chikmagalur)"
$ws.Range("G44").Value = "Shivamogga (Shimoga)"
$ws.Range("G47").Value = "Chikkamagaluru (Chikmagalur)"
$ws.Range("G48").Value = "Kalaburagi (Gulbarga)"
$ws.Range("G52").Value = "Shivamogga (Shimoga)"
$ws.Range("G53").Value = "Shivamogga (Shimoga)"
